$d = $word.ActiveDocument

# Locate the paragraph that currently reads "Escritura de archivos de texto"
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Escritura de archivos de texto*") {
        $target = $p
        break
    }
}
if ($target -eq $null) {
    throw "Could not locate target paragraph"
}

# Replace that paragraph (and refresh word/numbering.xml) via a flat-OPC
# WordprocessingML package: turns it into a bulleted list item, appends the
# new sentence about saving/reading the light/dark preference to a text
# file, and relocates the _GoBack bookmark to the end of the new text
# (removing it from its old spot at the very end of the document).
$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:wpc="http://schemas.microsoft.com/office/word/2010/wordprocessingCanvas" xmlns:cx="http://schemas.microsoft.com/office/drawing/2014/chartex" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w15="http://schemas.microsoft.com/office/word/2012/wordml" xmlns:w16se="http://schemas.microsoft.com/office/word/2015/wordml/symex" xmlns:wpg="http://schemas.microsoft.com/office/word/2010/wordprocessingGroup" xmlns:wpi="http://schemas.microsoft.com/office/word/2010/wordprocessingInk" xmlns:wne="http://schemas.microsoft.com/office/word/2006/wordml" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" mc:Ignorable="w14 w15 w16se wp14"><w:body><w:p w:rsidR="00986348" w:rsidRPr="00BC6C5E" w:rsidRDefault="00986348" w:rsidP="00BC6C5E"><w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:b/><w:u w:val="single"/></w:rPr></w:pPr><w:r w:rsidRPr="00BC6C5E"><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:b/><w:u w:val="single"/></w:rPr><w:t>Veterinaria Pequeñas Patitas</w:t></w:r></w:p><w:p w:rsidR="00986348" w:rsidRPr="00986348" w:rsidRDefault="00986348"><w:pPr><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/></w:rPr></w:pPr></w:p><w:p w:rsidR="00986348" w:rsidRPr="00BC6C5E" w:rsidRDefault="00986348"><w:pPr><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:b/></w:rPr></w:pPr><w:r w:rsidRPr="00BC6C5E"><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:b/></w:rPr><w:t xml:space="preserve">Corrección de errores del TP3 </w:t></w:r></w:p><w:p w:rsidR="00986348" w:rsidRDefault="00986348"><w:pPr><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:color w:val="FF0000"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/></w:rPr><w:t xml:space="preserve">Pruebas unitarias: se corrige la prueba unitaria </w:t></w:r><w:r w:rsidRPr="00986348"><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:color w:val="FF0000"/></w:rPr><w:t>TODO</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:color w:val="FF0000"/></w:rPr><w:t xml:space="preserve"> error en el constructor utiliza dos </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:color w:val="FF0000"/></w:rPr><w:t>parametros</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p w:rsidR="00986348" w:rsidRDefault="00986348"><w:pPr><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:color w:val="FF0000"/></w:rPr></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:color w:val="FF0000"/></w:rPr><w:t>Excepciones .</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:color w:val="FF0000"/></w:rPr><w:t xml:space="preserve"> excepción por defecto</w:t></w:r></w:p><w:p w:rsidR="00986348" w:rsidRDefault="00986348"><w:pPr><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:color w:val="FF0000"/></w:rPr><w:t xml:space="preserve">Explota cuando lee un </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:color w:val="FF0000"/></w:rPr><w:t>json</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:color w:val="FF0000"/></w:rPr><w:t xml:space="preserve"> que no existe porque la lista de paseos queda en </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:color w:val="FF0000"/></w:rPr><w:t>null</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:color w:val="FF0000"/></w:rPr><w:t xml:space="preserve"> en mostrar lista paseos</w:t></w:r></w:p><w:p w:rsidR="00986348" w:rsidRDefault="00986348"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/></w:rPr><w:t>Escritura de archivos de texto</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/></w:rPr><w:t>: Se guarda en un archivo de texto en el escritorio la actual preferencia de modo claro o modo oscuro. Se lee la preferencia desde allí también.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p w:rsidR="00986348" w:rsidRPr="00986348" w:rsidRDefault="00986348"><w:pPr><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/></w:rPr></w:pPr></w:p><w:p w:rsidR="00A2759A" w:rsidRPr="00BC6C5E" w:rsidRDefault="00986348"><w:pPr><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:b/></w:rPr></w:pPr><w:r w:rsidRPr="00BC6C5E"><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:b/></w:rPr><w:t>Temas vistos entre la clase 15 a 20</w:t></w:r></w:p><w:p w:rsidR="00986348" w:rsidRDefault="00986348" w:rsidP="00986348"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/></w:rPr></w:pPr><w:r w:rsidRPr="00986348"><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/></w:rPr><w:t>conexión</w:t></w:r><w:r w:rsidRPr="00986348"><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/></w:rPr><w:t xml:space="preserve"> a data base: Se </w:t></w:r><w:r w:rsidRPr="00986348"><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/></w:rPr><w:t>envían</w:t></w:r><w:r w:rsidRPr="00986348"><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/></w:rPr><w:t xml:space="preserve"> datos sobre la </w:t></w:r><w:r w:rsidRPr="00986348"><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/></w:rPr><w:t>vacunación</w:t></w:r><w:r w:rsidRPr="00986348"><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/></w:rPr><w:t xml:space="preserve"> canina y felina a la base de datos de un Instituto de Zoonosis. </w:t></w:r><w:r w:rsidRPr="00986348"><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/></w:rPr><w:t>También</w:t></w:r><w:r w:rsidRPr="00986348"><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/></w:rPr><w:t xml:space="preserve"> se recupera desde esa base de datos el total de vacunas aplicadas por la Veterinaria. Se adjunta el </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00986348"><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/></w:rPr><w:t>backup</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00986348"><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/></w:rPr><w:t xml:space="preserve"> de la base de datos en un archivo llamado </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00986348"><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/></w:rPr><w:t>INSTITUTOZOONOSIS.bak</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p w:rsidR="00986348" w:rsidRPr="00986348" w:rsidRDefault="00986348" w:rsidP="00986348"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/></w:rPr></w:pPr></w:p><w:p w:rsidR="00986348" w:rsidRDefault="00986348" w:rsidP="00986348"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Segoe UI Symbol" w:hAnsi="Segoe UI Symbol" w:cs="Segoe UI Symbol"/></w:rPr><w:t>D</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/></w:rPr><w:t>elegados y</w:t></w:r><w:r w:rsidRPr="00986348"><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00986348"><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/></w:rPr><w:t>multithreading</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00986348"><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/></w:rPr><w:t xml:space="preserve">: Hay una </w:t></w:r><w:r w:rsidRPr="00986348"><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/></w:rPr><w:t>función</w:t></w:r><w:r w:rsidRPr="00986348"><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/></w:rPr><w:t xml:space="preserve"> que simula una demora por carga de datos, que corre en un hilo diferente del hilo principal y utiliza delegados en la </w:t></w:r><w:r w:rsidRPr="00986348"><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/></w:rPr><w:t>invocación</w:t></w:r><w:r w:rsidRPr="00986348"><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/></w:rPr><w:t xml:space="preserve"> de los </w:t></w:r><w:r w:rsidRPr="00986348"><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/></w:rPr><w:t>métodos</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/></w:rPr><w:t>.</w:t></w:r></w:p><w:p w:rsidR="003620DD" w:rsidRPr="003620DD" w:rsidRDefault="003620DD" w:rsidP="003620DD"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/></w:rPr></w:pPr></w:p><w:p w:rsidR="003620DD" w:rsidRDefault="003620DD" w:rsidP="00986348"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/></w:rPr><w:t>También se usa un delegado en el evento.</w:t></w:r></w:p><w:p w:rsidR="00986348" w:rsidRPr="00986348" w:rsidRDefault="00986348" w:rsidP="00986348"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/></w:rPr></w:pPr></w:p><w:p w:rsidR="00986348" w:rsidRPr="00986348" w:rsidRDefault="00986348" w:rsidP="00986348"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/></w:rPr></w:pPr></w:p><w:p w:rsidR="00986348" w:rsidRDefault="00986348" w:rsidP="00986348"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Segoe UI Symbol" w:hAnsi="Segoe UI Symbol" w:cs="Segoe UI Symbol"/></w:rPr><w:t>E</w:t></w:r><w:r w:rsidRPr="00986348"><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/></w:rPr><w:t>xpresiones lambda: se utilizan en varios puntos</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/></w:rPr><w:t>,</w:t></w:r><w:r w:rsidRPr="00986348"><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/></w:rPr><w:t xml:space="preserve"> por </w:t></w:r><w:proofErr w:type="gramStart"/><w:r w:rsidRPr="00986348"><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/></w:rPr><w:t>ejemplo</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r w:rsidRPr="00986348"><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/></w:rPr><w:t xml:space="preserve"> en las propiedades de la clase Perro y de la clase Animal</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/></w:rPr><w:t xml:space="preserve">. </w:t></w:r><w:r w:rsidRPr="00986348"><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/><w:color w:val="FF0000"/></w:rPr><w:t>TODO</w:t></w:r></w:p><w:p w:rsidR="00986348" w:rsidRPr="00986348" w:rsidRDefault="00986348" w:rsidP="00986348"><w:pPr><w:ind w:left="360"/><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/></w:rPr></w:pPr></w:p><w:p w:rsidR="00986348" w:rsidRPr="00986348" w:rsidRDefault="00986348" w:rsidP="00986348"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Segoe UI Symbol" w:hAnsi="Segoe UI Symbol" w:cs="Segoe UI Symbol"/></w:rPr><w:t>E</w:t></w:r><w:r w:rsidRPr="00986348"><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/></w:rPr><w:t>ventos</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/></w:rPr><w:t>:</w:t></w:r><w:r w:rsidR="005546F0"><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/></w:rPr><w:t xml:space="preserve"> La app tiene un modo claro y un modo oscuro, según se elija uno u otro</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/></w:rPr><w:t xml:space="preserve"> se dispara un evento</w:t></w:r><w:r w:rsidR="005546F0"><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana"/></w:rPr><w:t xml:space="preserve"> que realiza el cambio.</w:t></w:r></w:p><w:sectPr w:rsidR="00986348" w:rsidRPr="00986348"><w:pgSz w:w="11906" w:h="16838"/><w:pgMar w:top="1417" w:right="1701" w:bottom="1417" w:left="1701" w:header="708" w:footer="708" w:gutter="0"/><w:cols w:space="708"/><w:docGrid w:linePitch="360"/></w:sectPr></w:body></w:document>
</pkg:xmlData>
</pkg:part>
<pkg:part pkg:name="/word/numbering.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.numbering+xml">
<pkg:xmlData>
<w:numbering xmlns:wpc="http://schemas.microsoft.com/office/word/2010/wordprocessingCanvas" xmlns:cx="http://schemas.microsoft.com/office/drawing/2014/chartex" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w15="http://schemas.microsoft.com/office/word/2012/wordml" xmlns:w16se="http://schemas.microsoft.com/office/word/2015/wordml/symex" xmlns:wpg="http://schemas.microsoft.com/office/word/2010/wordprocessingGroup" xmlns:wpi="http://schemas.microsoft.com/office/word/2010/wordprocessingInk" xmlns:wne="http://schemas.microsoft.com/office/word/2006/wordml" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" mc:Ignorable="w14 w15 w16se wp14"><w:abstractNum w:abstractNumId="0" w15:restartNumberingAfterBreak="0"><w:nsid w:val="36813BAD"/><w:multiLevelType w:val="hybridMultilevel"/><w:tmpl w:val="5D02A7F4"/><w:lvl w:ilvl="0" w:tplc="2C0A000D"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="720" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Wingdings" w:hAnsi="Wingdings" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="1" w:tplc="2C0A0003" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val="o"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="1440" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="2" w:tplc="2C0A0005" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="2160" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Wingdings" w:hAnsi="Wingdings" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="3" w:tplc="2C0A0001" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="2880" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Symbol" w:hAnsi="Symbol" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="4" w:tplc="2C0A0003" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val="o"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="3600" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="5" w:tplc="2C0A0005" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="4320" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Wingdings" w:hAnsi="Wingdings" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="6" w:tplc="2C0A0001" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="5040" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Symbol" w:hAnsi="Symbol" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="7" w:tplc="2C0A0003" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val="o"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="5760" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="8" w:tplc="2C0A0005" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="6480" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Wingdings" w:hAnsi="Wingdings" w:hint="default"/></w:rPr></w:lvl></w:abstractNum><w:abstractNum w:abstractNumId="1" w15:restartNumberingAfterBreak="0"><w:nsid w:val="7A8C19A2"/><w:multiLevelType w:val="hybridMultilevel"/><w:tmpl w:val="303E4AB4"/><w:lvl w:ilvl="0" w:tplc="2C0A000D"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="720" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Wingdings" w:hAnsi="Wingdings" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="1" w:tplc="2C0A0003" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val="o"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="1440" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="2" w:tplc="2C0A0005" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="2160" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Wingdings" w:hAnsi="Wingdings" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="3" w:tplc="2C0A0001" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="2880" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Symbol" w:hAnsi="Symbol" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="4" w:tplc="2C0A0003" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val="o"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="3600" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="5" w:tplc="2C0A0005" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="4320" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Wingdings" w:hAnsi="Wingdings" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="6" w:tplc="2C0A0001" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="5040" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Symbol" w:hAnsi="Symbol" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="7" w:tplc="2C0A0003" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val="o"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="5760" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="8" w:tplc="2C0A0005" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val=""/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="6480" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Wingdings" w:hAnsi="Wingdings" w:hint="default"/></w:rPr></w:lvl></w:abstractNum><w:num w:numId="1"><w:abstractNumId w:val="1"/></w:num><w:num w:numId="2"><w:abstractNumId w:val="0"/></w:num></w:numbering>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$target.Range.InsertXML($xml)

Write-Output "done"
